$wb = $excel.ActiveWorkbook

# Add the new "low" sheet after "exist" (so the tab order becomes new, exist, low)
$existSheet = $wb.Worksheets.Item("exist")
[void]$existSheet.Range("A1:E12").Select()
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $existSheet)
$ws.Name = "low"

# Headers (same as the other two sheets)
$ws.Range("A1").Value = "REGION_GEA"
$ws.Range("B1").Value = "RegNum"
$ws.Range("C1").Value = "arch"
$ws.Range("D1").Value = "rural"
$ws.Range("E1").Value = "urban"

# Region names / numbers, identical across all three sheets
$regions = @("AFR","CPA","EEU","FSU","LAC","MEA","NAM","PAO","PAS","SAS","WEU")
for ($i = 0; $i -lt $regions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $regions[$i]
    $ws.Cells.Item($row, 2).Value = $i + 1
}

# Column C: literal "low" in row 2, then formulas referencing the row above.
# Row 3 gets its own (non-shared) formula, rows 4:12 are filled as one block
# so Excel groups them into a single shared-formula range (matches the
# pattern already used by columns D/E on every sheet in this workbook).
$ws.Range("C2").Value = "low"
$ws.Range("C3").Formula = "=C2"
$ws.Range("C4:C12").Formula = "=C3"

# Columns D/E: literal seed values in row 2, then formulas adding 2 each row.
# D3/E3 are set individually (plain, non-shared formulas) and D4:E12 is
# filled as one block (shared-formula group) -- matching the layout already
# used on the "new" and "exist" sheets.
$ws.Range("D2").Value = 44
$ws.Range("E2").Value = 45
$ws.Range("D3").Formula = "=D2+2"
$ws.Range("E3").Formula = "=E2+2"
$ws.Range("D4:E12").Formula = "=D3+2"

$ws.Calculate()

# Make "low" the active (selected) tab, matching the saved file's activeTab
[void]$ws.Select()
[void]$ws.Range("H46").Select()
